$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Fix the font-name typo across the whole document: TimesNewToman -> Times New Roman
# ---------------------------------------------------------------------------
$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"

# Helper pattern: use Find/Replace (wdReplaceOne) for simple 1:1 text swaps.
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)

# ---------------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Unraveling the Complexities of Surgical Robotics: Precision, Innovation, and Challenges",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The Art of Persuasion: Understanding and Applying Rhetorical Devices in Public Speaking", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Author name ("Dr. Isabella Phillips" spans 3 runs -> single run "Robert Hill")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Dr. Isabella Phillips",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Robert Hill", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Email address (local-part run + domain run; the middle "." run is kept)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "isabellaphillips@premiumresearch",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "rhill@school", 2) | Out-Null

$d.Content.Find.Execute(
    "com",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "edu", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Body paragraph - sentence-by-sentence replacements
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "In the realm of modern healthcare, surgical robotics has emerged as a transformative force, redefining the boundaries of precision and patient outcomes",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Oratory, the art of effective public speaking, is a skill that can empower individuals to convey their thoughts, influence audiences, and shape opinions", 2) | Out-Null

$d.Content.Find.Execute(
    " As a technology that seamlessly intertwines engineering prowess with medical expertise, surgical robotics unveils a future where minimally invasive procedures, unparalleled accuracy, and expedited recovery times coalesce",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " Skilled speakers harness the power of words and rhetorical devices to create compelling messages that persuade and inspire", 2) | Out-Null

$d.Content.Find.Execute(
    " Journey with us as we delve into the intricacies of surgical robotics, illuminating its groundbreaking applications while acknowledging the complexities and challenges that accompany this revolutionary technology",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " In this essay, we will explore the world of persuasion, unveiling the techniques and strategies speakers use to captivate and sway their audiences", 2) | Out-Null

# New sentence pair inserted right before the first <w:br/><w:br/> (after "...audiences.")
$r = $d.Content
$r.Find.Execute(" In this essay, we will explore the world of persuasion, unveiling the techniques and strategies speakers use to captivate and sway their audiences", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$afterPeriod = $d.Range($r.End, $r.End + 1)
$insPos = $afterPeriod.End
$ins = $d.Range($insPos, $insPos)
$ins.InsertAfter(" We will delve into the nuances of language, the impact of effective delivery, and the art of connecting with listeners on an emotional level.")

# "Within the operating room..." -> "Persuasion is an intricate dance between speaker, audience, and context"
$d.Content.Find.Execute(
    "Within the operating room, surgical robots extend the reach of the surgeon, enabling them to navigate intricate anatomical structures with a finesse that surpasses human capabilities",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Persuasion is an intricate dance between speaker, audience, and context", 2) | Out-Null

$d.Content.Find.Execute(
    " The fusion of computer-aided visualization with robotic dexterity empowers surgeons to perform intricate tasks with enhanced precision, minimizing tissue trauma and reducing surgical complications",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " Understanding the dynamics of this relationship is key to crafting persuasive messages", 2) | Out-Null

$d.Content.Find.Execute(
    " Furthermore, surgical robotics enables remote procedures, connecting surgeons and patients across vast distances, potentially revolutionizing access to specialized surgical care",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " Additionally, knowing your audience, their needs, and their values allows you to tailor your speech to resonate with them", 2) | Out-Null

# New sentence pair inserted right before the second <w:br/><w:br/> (after "...resonate with them.")
$r2 = $d.Content
$r2.Find.Execute(" Additionally, knowing your audience, their needs, and their values allows you to tailor your speech to resonate with them", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(0)
$afterPeriod2 = $d.Range($r2.End, $r2.End + 1)
$insPos2 = $afterPeriod2.End
$ins2 = $d.Range($insPos2, $insPos2)
$ins2.InsertAfter(" Whether trying to win over a vote, promote a cause, or simply share an idea, persuasive speaking is an invaluable skill that can be learned and refined.")

$d.Content.Find.Execute(
    "The rise of surgical robotics has not been without its share of challenges",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "In this multifaceted journey, we will uncover the secrets of effective persuasion, unravelling the mysteries of rhetoric and exploring the psychology of influence", 2) | Out-Null

$d.Content.Find.Execute(
    " These marvels of engineering are undeniably complex, mandating extensive training for surgeons to master their operation",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " By gaining these insights, we can become more effective communicators, leaving lasting impressions on our listeners and leaving impacts that can change the world", 2) | Out-Null

# Remove the old trailing sentences ("Moreover ... the capabilities of the robot.") entirely.
$r3 = $d.Content
$r3.Find.Execute(" By gaining these insights, we can become more effective communicators, leaving lasting impressions on our listeners and leaving impacts that can change the world", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Collapse(0)
$startDel = $d.Range($r3.End, $r3.End + 1).End   # past the "." that follows

$r4 = $d.Content
$r4.Find.Execute(" The symbiosis between humans and machines remains a critical area of exploration, as surgeons strive to strike an optimal balance between the expertise of the surgeon and the capabilities of the robot", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Collapse(0)
$endDel = $d.Range($r4.End, $r4.End + 1).End    # past the final "."

$d.Range($startDel, $endDel).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 6. Summary paragraph - sentence-by-sentence replacements
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Surgical robotics has indubitably ushered in a new era of surgical precision, paving the way for minimally invasive procedures and improved patient outcomes",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This essay explores the art of persuasion, unveiling the techniques and strategies speakers use to captivate and sway their audiences", 2) | Out-Null

$d.Content.Find.Execute(
    " However, the complexities of this technology, coupled with financial implications and the dynamic relationship between humans and machines, necessitate ongoing research and advancement",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " We delve into the nuances of language, the impact of effective delivery, and the art of connecting with listeners on an emotional level", 2) | Out-Null

$d.Content.Find.Execute(
    " As surgical robotics continues to evolve, it holds the potential to transform healthcare landscapes globally, enhancing accessibility to specialized surgical care and empowering surgeons with unprecedented precision",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " By understanding the dynamics of persuasion and knowing your audience, you can craft persuasive messages that resonate with them", 2) | Out-Null

# New sentence pair appended at the end of the Summary paragraph
$r5 = $d.Content
$r5.Find.Execute(" By understanding the dynamics of persuasion and knowing your audience, you can craft persuasive messages that resonate with them", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r5.Collapse(0)
$afterPeriod5 = $d.Range($r5.End, $r5.End + 1)
$insPos5 = $afterPeriod5.End
$ins5 = $d.Range($insPos5, $insPos5)
$ins5.InsertAfter(" Unlock the secrets of effective persuasion, unravel the mysteries of rhetoric, and explore the psychology of influence to become a more effective communicator and leave lasting impressions on your listeners.")

# ---------------------------------------------------------------------------
# 7. Add a new empty paragraph at the very end of the document (before sectPr)
# ---------------------------------------------------------------------------
$endOfDoc = $d.Content.End
$d.Range($endOfDoc, $endOfDoc).InsertParagraphAfter() | Out-Null

Write-Host "Edit complete."
